$wb = $excel.ActiveWorkbook

# 1. Move the "Items" sheet so it sits right after "Locations"
#    (new order: Races, Classes, Players, NPCs, Locations, Items,
#     Time-Pieces, Time-Anomalies, Quests, Snitel)
$itemsSheet = $wb.Worksheets.Item("Items")
$locationsSheet = $wb.Worksheets.Item("Locations")
$itemsSheet.Move($null, $locationsSheet)

# 2. Update the Items sheet headers / data
$ws = $wb.Worksheets.Item("Items")

$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Attributes"
$ws.Range("D1").Value = "Rarity"
$ws.Range("E1").Value = "Owned By"
$ws.Range("A1:E1").Font.Bold = $true

$ws.Range("A2").Value = "Sword"
$ws.Range("B2").Value = "A sword"
$ws.Range("C2").Value = "'+2 attack"
$ws.Range("D2").Value = "Common"
$ws.Range("E2").Value = "None"

# 3. Make the Items sheet the active tab, with A2 selected
$ws.Activate()
$ws.Range("A2").Select()
